# Add a new "15-10-2020" column (AD) to the deceased-cases time-series sheet.
# This mirrors the existing AC column: same header style/border/bold and
# the day's cumulative deceased-case counts for each state/UT row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell AD1 -------------------------------------------------
# Copy the formatting (bold, centered, thin border) of the adjacent AC1
# header cell onto AD1, then set its text value.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = "15-10-2020"

# --- Data cells AD2:AD36 (cumulative deceased counts for 15-10-2020) --
$values = @{
    2  = 55
    3  = 6319
    4  = 29
    5  = 834
    6  = 967
    7  = 199
    8  = 1339
    9  = 2
    10 = 5898
    11 = 519
    12 = 3595
    13 = 1614
    14 = 255
    15 = 1352
    16 = 811
    17 = 10198
    18 = 1066
    19 = 64
    20 = 2686
    21 = 40859
    22 = 103
    23 = 70
    24 = 0
    25 = 22
    26 = 1062
    27 = 568
    28 = 3925
    29 = 1694
    30 = 59
    31 = 10423
    32 = 1249
    33 = 319
    34 = 796
    35 = 6507
    36 = 5808
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 30).Value = $values[$row]
}
